# The workbook currently has a single worksheet named "ODI Batting"
# (sheetId 1). The target state has two sheets:
#   1. "Player Info"  (sheetId 1, reuses the original sheet's identity)
#   2. "ODI Batting"  (sheetId 2, a new sheet holding the original batting
#                       data, with the MATCH_CARD_LINK column turned into a
#                       MATCH_CODE column)
#
# To reproduce that with the original sheet's identity (rId1/sheetId1)
# ending up on "Player Info", we:
#   - create a new sheet right after the existing one,
#   - copy the existing batting data into it and tweak column D,
#   - rename the new sheet to "ODI Batting",
#   - clear + rename the original sheet to "Player Info" and populate it.

$wb = $excel.ActiveWorkbook

$odi = $wb.Worksheets.Item("ODI Batting")

# --- 1) Create the new sheet that will carry the batting data forward ---
$newOdi = $wb.Worksheets.Add($null, $odi)
$newOdi.Name = "ODI Batting (new)"

# Duplicate the original sheet's data/formatting onto the new sheet.
$odi.Cells.Copy($newOdi.Range("A1"))

# Apply the requested column-D change: MATCH_CARD_LINK -> MATCH_CODE,
# and the URL value -> the bare match code "4705" (stored as text).
$newOdi.Range("D1").Value = "MATCH_CODE"
$newOdi.Range("D2").NumberFormat = "@"
$newOdi.Range("D2").Value = "4705"

# --- 2) Turn the original sheet into the new "Player Info" sheet ---
$odi.Cells.Clear()
$odi.Name = "Player Info"

$odi.Range("A1").Value = "ID"
$odi.Range("B1").Value = "NAME"
$odi.Range("C1").Value = "BATTING_HAND"
$odi.Range("D1").Value = "BOWL_STYLE"

# Store the data row as text (matching the source data's string typing).
$odi.Range("A2:D2").NumberFormat = "@"
$odi.Range("A2").Value = "7130"
$odi.Range("B2").Value = "Liam Robert Naylor"
$odi.Range("C2").Value = "Right Handed"
$odi.Range("D2").Value = "Right Arm Medium"

# Bold, centered, bordered header row - matching the style already used
# for headers on the batting sheet.
$header = $odi.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108  # xlCenter
$header.VerticalAlignment = -4160    # xlTop
$header.Borders.LineStyle = 1        # xlContinuous

# --- 3) Finish renaming / ordering ---
$newOdi.Name = "ODI Batting"

# "Player Info" is first and stays the active tab (activeTab stays 0).
$odi.Activate()
